# daily auto push: 2026-01-19 02:34 UTC
# A new sample row for 2026/01/19 07:00 (ranking 20) was inserted into the
# daily log right before the existing "2026/12/29" block, pushing every
# row from the old r654 onward down by one (old r695 -> new r696).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 654; everything currently at/after 654
# (the 2026/12/29 ... 2027/01/05 block) shifts down one row.
$ws.Rows.Item(654).Insert()

# Column A holds dates formatted as plain text (e.g. "2026/12/29"), not
# real Excel date serials, so force text formatting before writing the
# slash-separated value - otherwise Excel would coerce it to a date.
$ws.Range("A654").NumberFormat = "@"
$ws.Range("A654").Value = "2026/01/19"
$ws.Range("A654").Style = "Normal"

$ws.Range("B654").Value = "月"
$ws.Range("C654").Value = 7
$ws.Range("D654").Value = 20
